$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 99
$ws.Range("I4").Value = 99
$ws.Range("K4").Value = 99
$ws.Range("M4").Value = 15
$ws.Range("H33").Value = 979.6667
$ws.Range("I33").Value = 1145.25
$ws.Range("K33").Value = 1145.25
$ws.Range("M33").Value = -916.25
$ws.Range("H88").Value = 2624.5
$ws.Range("I88").Value = 1833
$ws.Range("K88").Value = 1833
$ws.Range("M88").Value = -1427
$ws.Range("H91").Value = 2624.5
$ws.Range("I91").Value = 1833
$ws.Range("K91").Value = 1833
$ws.Range("M91").Value = -429
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("N116").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 450.44446
$ws.Range("I4").Value = 576
$ws.Range("J4").Value = 199.33333
$ws.Range("K4").Value = 576
$ws.Range("L4").Value = 199.33333
$ws.Range("M4").Value = -460
$ws.Range("N4").Value = -431.33333
$ws.Range("H45").Value = 2633.3333
$ws.Range("I45").Value = 2633.3333
$ws.Range("K45").Value = 2633.3333
$ws.Range("M45").Value = -2256.3333
$ws.Range("H59").Value = 47000
$ws.Range("J59").Value = 47000
$ws.Range("L59").Value = 47000
$ws.Range("N59").Value = -48608
$ws.Range("H96").Value = 35663.43
$ws.Range("J96").Value = 35663.43
$ws.Range("L96").Value = 35663.43
$ws.Range("N96").Value = -41155.43
$ws.Range("H102").Value = 1664
$ws.Range("I102").Value = 1358
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 1358
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 264
$ws.Range("N102").Value = -6744

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 982.8333
$ws.Range("I22").Value = 899.25
$ws.Range("K22").Value = 899.25
$ws.Range("M22").Value = -726.25
$ws.Range("H68").Value = 33000
$ws.Range("J68").Value = 33000
$ws.Range("L68").Value = 33000
$ws.Range("N68").Value = -34622
$ws.Range("H71").Value = 33000
$ws.Range("J71").Value = 33000
$ws.Range("L71").Value = 99000
$ws.Range("N71").Value = -107112
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H105").Value = 2711.25
$ws.Range("I105").Value = 2711.25
$ws.Range("K105").Value = 2711.25
$ws.Range("M105").Value = -964.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 298.1111
$ws.Range("I7").Value = 271.6
$ws.Range("J7").Value = 331.25
$ws.Range("K7").Value = 271.6
$ws.Range("L7").Value = 331.25
$ws.Range("M7").Value = -158.6
$ws.Range("N7").Value = -557.25
$ws.Range("H16").Value = 1399.1428
$ws.Range("I16").Value = 1099
$ws.Range("K16").Value = 1099
$ws.Range("M16").Value = -812
$ws.Range("H22").Value = 952.2
$ws.Range("J22").Value = 940.5
$ws.Range("L22").Value = 940.5
$ws.Range("N22").Value = -1640.5
$ws.Range("H33").Value = 19553.666
$ws.Range("I33").Value = 11464.4
$ws.Range("J33").Value = 60000
$ws.Range("K33").Value = 11464.4
$ws.Range("L33").Value = 60000
$ws.Range("M33").Value = -11085.4
$ws.Range("N33").Value = -60758
$ws.Range("H39").Value = 50000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 50000
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -50782
$ws.Range("H49").Value = 50000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 50000
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -50364
$ws.Range("H113").Value = 1399.1428
$ws.Range("I113").Value = 1099
$ws.Range("K113").Value = 1099
$ws.Range("M113").Value = 1071

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2773.7273
$ws.Range("I14").Value = 2773.7273
$ws.Range("K14").Value = 8321.1819
$ws.Range("M14").Value = -8148.1819
$ws.Range("H50").Value = 1476
$ws.Range("I50").Value = 452
$ws.Range("J50").Value = 2500
$ws.Range("K50").Value = 1356
$ws.Range("L50").Value = 7500
$ws.Range("M50").Value = -875
$ws.Range("N50").Value = -8462
$ws.Range("H53").Value = 1476
$ws.Range("I53").Value = 452
$ws.Range("J53").Value = 2500
$ws.Range("K53").Value = 1356
$ws.Range("L53").Value = 7500
$ws.Range("M53").Value = -875
$ws.Range("N53").Value = -8462

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 311.23077
$ws.Range("I2").Value = 472.125
$ws.Range("J2").Value = 53.8
$ws.Range("K2").Value = 472.125
$ws.Range("L2").Value = 53.8
$ws.Range("M2").Value = -359.125
$ws.Range("N2").Value = -279.8
$ws.Range("H11").Value = 286428.56
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = -999861
$ws.Range("N11").Value = -1278
$ws.Range("H80").Value = 2994.25
$ws.Range("I80").Value = 2825.6667
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2825.6667
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1827.6667
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 2994.25
$ws.Range("I83").Value = 2825.6667
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 14128.3335
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -9136.333500000001
$ws.Range("N83").Value = -27484

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1599.1538
$ws.Range("I22").Value = 1818.9
$ws.Range("J22").Value = 866.6667
$ws.Range("K22").Value = 1818.9
$ws.Range("L22").Value = 866.6667
$ws.Range("M22").Value = -1523.9
$ws.Range("N22").Value = -1456.6667
$ws.Range("H27").Value = 1599.1538
$ws.Range("I27").Value = 1818.9
$ws.Range("J27").Value = 866.6667
$ws.Range("K27").Value = 1818.9
$ws.Range("L27").Value = 866.6667
$ws.Range("M27").Value = -1711.9
$ws.Range("N27").Value = -1080.6667
$ws.Range("H46").Value = 2356.4285
$ws.Range("I46").Value = 699
$ws.Range("K46").Value = 699
$ws.Range("M46").Value = -511

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2220.8
$ws.Range("I136").Value = 1918
$ws.Range("J136").Value = 2675
$ws.Range("K136").Value = 5754
$ws.Range("L136").Value = 8025
$ws.Range("M136").Value = -3204
$ws.Range("N136").Value = -13125
